$d = $word.ActiveDocument

# 1. Delete the leading empty paragraph (first paragraph in the document,
#    which holds no text, just an <w:rtl w:val="0"/> run property).
$d.Paragraphs(1).Range.Delete()

# 2. "de la dicte crouste" -> "de ladicte crouste"
$d.Content.Find.Execute("de la dicte crouste", $true, $false, $false, $false, $false,
                         $true, 1, $false, "de ladicte crouste", 2)

# 3. "quattres petits" -> "quattre petits"
$d.Content.Find.Execute("quattres petits", $true, $false, $false, $false, $false,
                         $true, 1, $false, "quattre petits", 2)
